$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product ("Tubos para armar cigarrillos Stamps") was added to the
# price list as row 5, pushing the existing rows (old 5..36) down to 6..37.
$ws.Rows.Item(5).Insert()

# The inserted row should keep the same "ImagenExactaDelArticulo" cell
# formatting (style) used by the rest of the data rows, so copy that
# formatting down from the row that now sits right below it.
$ws.Range("O6").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A5").Value = 7798205440198
$ws.Range("B5").Value = "Tubos"
$ws.Range("C5").Value = "para armar"
$ws.Range("D5").Value = "cigarrillos"
$ws.Range("E5").Value = "Stamps"
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = "und."
$ws.Range("H5").Value = "Caja"
$ws.Range("I5").Value = "Tabaco"
$ws.Range("J5").Value = "Argentina"
$ws.Range("K5").Value = 6
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = $false
$ws.Range("N5").Value = "C:\VentaSoft\Imágenes de artículos\7798205440198.png"
$ws.Range("O5").Value = $true
